# edit.ps1 - applies the "added hint on addr2line and nm usage" commit
#
# Summary of changes:
#  1. Slide 2 ("Overview"): merge 3 runs of the 1st content paragraph into one run
#     "Memory handling related crashes".
#  2. Slide 4 ("How Valgrind works..."): merge the trailing "of " + "concurrency."
#     runs into a single "of concurrency." run (leaving "because " untouched).
#  3. Slide 5 ("Crash investigator - idea behind"):
#       a. merge the 2 title runs into a single "Crash investigator - idea behind" run.
#       b. merge the 3 runs of the 1st content paragraph into one run.
#  4. Slide 6 ("Some cases of indirect double/free"): merge the trailing
#     "... proper " + "synchronization." runs into one run.
#  5. Append a new slide (8th / last) re-using the "Title and Content" layout,
#     with an empty title and a content placeholder containing the addr2line/nm
#     hint text.

$p = $ppt.ActivePresentation

function Merge-ParagraphText($para, $finalText) {
    # Forcing a text change (even transiently) makes the host collapse a
    # paragraph's multiple runs down to a single run that carries the
    # formatting of what used to be its first run - exactly matching how
    # PowerPoint coalesces runs when a user retypes/merges text manually.
    $para.Text = "##TMP##"
    $para.Text = $finalText
}

# ---------------------------------------------------------------------------
# 1. Slide 2 - "Overview": "Memory " + "handling related " + "crashes"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
Merge-ParagraphText $tr2.Paragraphs(1, 1) "Memory handling related crashes"

# ---------------------------------------------------------------------------
# 2. Slide 4 - "How Valgrind works...": "...because " + "of " + "concurrency."
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$full4 = $tr4.Text
$idx4 = $full4.IndexOf("of concurrency.")
$sub4 = $tr4.Characters($idx4 + 1, 15)
Merge-ParagraphText $sub4 "of concurrency."

# ---------------------------------------------------------------------------
# 3. Slide 5 - "Crash investigator - idea behind"
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)

# 3a. Title: "Crash investigator - " + "idea behind"
$titleTr5 = $s5.Shapes.Item(1).TextFrame.TextRange
Merge-ParagraphText $titleTr5 "Crash investigator - idea behind"

# 3b. Content paragraph 1: "In order to trap ... following " + "can be " + "done"
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange
Merge-ParagraphText $tr5.Paragraphs(1, 1) "In order to trap memory allocations functions and make some analyze the following can be done"

# ---------------------------------------------------------------------------
# 4. Slide 6 - "Some cases of indirect double/free": "...proper " + "synchronization."
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange
Merge-ParagraphText $tr6.Paragraphs(3, 1) "Global buffers allocation/deallocation without proper synchronization."

# ---------------------------------------------------------------------------
# 5. New slide 8, appended at the end, "Title and Content" layout.
# ---------------------------------------------------------------------------
$layout = $p.SlideMaster.CustomLayouts.Item(2)
$newSlide = $p.Slides.AddSlide($p.Slides.Count + 1, $layout)

$newSlide.Shapes.Item(1).Name = "Title 1"
$newSlide.Shapes.Item(2).Name = "Content Placeholder 2"

# Title stays empty, but (as in the authored deck) carries a German
# language mark.
$titleTrNew = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleTrNew.Text = ""
$titleTrNew.LanguageID = "de-DE"

$line1 = "addr2line -e /home/kalantar/dev/crash_investigator/sys/focal/Debug/lib/libcrash_investigator_new_malloc_0020.so -f -C 0x7d75"
$line2 = "nm /home/kalantar/dev/crash_investigator/sys/focal/Debug/lib/libcrash_investigator_new_malloc_0020.so | grep _ZN18crash_investigator11CMemoryItem4InitEmNS_11FailureTypeEPvPNS_9BacktraceE"

$bodyTrNew = $newSlide.Shapes.Item(2).TextFrame.TextRange
$bodyTrNew.Text = $line1
$bodyTrNew.LanguageID = "de-DE"
$null = $bodyTrNew.InsertAfter("`r" + $line2)

# Re-fetch the range and stamp the language onto the 2nd paragraph too -
# inserting a new paragraph after setting LanguageID on the 1st does not
# automatically propagate it.
$bodyTrNew2 = $newSlide.Shapes.Item(2).TextFrame.TextRange
$bodyTrNew2.Paragraphs(2, 1).LanguageID = "de-DE"
